$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I5").Value = 0.5
$ws.Range("H8").Value = 0.5
$ws.Range("H9").Value = 0.5
$ws.Range("H10").Value = 0.5
$ws.Range("G11").Value = 5.0
$ws.Range("F14").Value = ""
